$wb = $excel.ActiveWorkbook

# Update the shared string "Ready for handoff" -> "In Translation" everywhere it appears
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns are E (zh-cn) and F (de-de), rows 2-4
foreach ($r in 2..4) {
    if ($overview.Cells.Item($r, 5).Value2 -eq "Ready for handoff") {
        $overview.Cells.Item($r, 5).Value = "In Translation"
    }
    if ($overview.Cells.Item($r, 6).Value2 -eq "Ready for handoff") {
        $overview.Cells.Item($r, 6).Value = "In Translation"
    }
}

# zh-cn / de-de sheets: Status column is C, rows 2-4
foreach ($r in 2..4) {
    if ($zhcn.Cells.Item($r, 3).Value2 -eq "Ready for handoff") {
        $zhcn.Cells.Item($r, 3).Value = "In Translation"
    }
    if ($dede.Cells.Item($r, 3).Value2 -eq "Ready for handoff") {
        $dede.Cells.Item($r, 3).Value = "In Translation"
    }
}

# The shorter replacement text narrows the status columns; set the new
# (re-measured) column widths directly to match the regenerated report.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
